$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2532085.8
$ws.Range("I18").Value = 3970848
$ws.Range("J18").Value = 14251.5
$ws.Range("K18").Value = 3970848
$ws.Range("L18").Value = 14251.5
$ws.Range("M18").Value = -3970564
$ws.Range("N18").Value = -14819.5
$ws.Range("H21").Value = 70019
$ws.Range("I21").Value = 70019
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 70019
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -69551
$ws.Range("H23").Value = 70019
$ws.Range("I23").Value = 70019
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 70019
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -69785
$ws.Range("H69").Value = 3668.6
$ws.Range("I69").Value = 3806.5
$ws.Range("J69").Value = 3576.6667
$ws.Range("K69").Value = 11419.5
$ws.Range("L69").Value = 10730.0001
$ws.Range("M69").Value = -10545.5
$ws.Range("N69").Value = -12478.0001
$ws.Range("H72").Value = 3668.6
$ws.Range("I72").Value = 3806.5
$ws.Range("J72").Value = 3576.6667
$ws.Range("K72").Value = 34258.5
$ws.Range("L72").Value = 32190.0003
$ws.Range("M72").Value = -29890.5
$ws.Range("N72").Value = -40926.0003
$ws.Range("H133").Value = 29513.334
$ws.Range("J133").Value = 29513.334
$ws.Range("L133").Value = 29513.334
$ws.Range("N133").Value = -39633.334
$ws.Range("H137").Value = 5993.0713
$ws.Range("I137").Value = 13600
$ws.Range("J137").Value = 2950.3
$ws.Range("K137").Value = 40800
$ws.Range("L137").Value = 8850.900000000001
$ws.Range("M137").Value = -38250
$ws.Range("N137").Value = -13950.9

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 400.2
$ws.Range("I5").Value = 300.25
$ws.Range("J5").Value = 800
$ws.Range("K5").Value = 300.25
$ws.Range("L5").Value = 800
$ws.Range("M5").Value = -188.25
$ws.Range("N5").Value = -1024
$ws.Range("H61").Value = 3698.0715
$ws.Range("I61").Value = 1666.6666
$ws.Range("J61").Value = 3941.84
$ws.Range("K61").Value = 1666.6666
$ws.Range("L61").Value = 3941.84
$ws.Range("M61").Value = -1454.6666
$ws.Range("N61").Value = -4365.84
$ws.Range("H136").Value = 3698.0715
$ws.Range("I136").Value = 1666.6666
$ws.Range("J136").Value = 3941.84
$ws.Range("K136").Value = 4999.9998
$ws.Range("L136").Value = 11825.52
$ws.Range("M136").Value = -2449.9998
$ws.Range("N136").Value = -16925.52

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 400.2
$ws.Range("I4").Value = 300.25
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 300.25
$ws.Range("L4").Value = 800
$ws.Range("M4").Value = -185.25
$ws.Range("N4").Value = -1030
$ws.Range("H19").Value = 58006
$ws.Range("J19").Value = 58006
$ws.Range("L19").Value = 58006
$ws.Range("N19").Value = -58352
$ws.Range("H97").Value = 14766.182
$ws.Range("I97").Value = 8485.6
$ws.Range("K97").Value = 8485.6
$ws.Range("M97").Value = -7494.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 114.44444
$ws.Range("I7").Value = 75
$ws.Range("J7").Value = 146
$ws.Range("K7").Value = 75
$ws.Range("L7").Value = 146
$ws.Range("M7").Value = 38
$ws.Range("N7").Value = -372
$ws.Range("H36").Value = 70052.75
$ws.Range("I36").Value = 70052
$ws.Range("J36").Value = 70053
$ws.Range("K36").Value = 70052
$ws.Range("L36").Value = 70053
$ws.Range("M36").Value = -69664
$ws.Range("N36").Value = -70829
$ws.Range("H40").Value = 70052.75
$ws.Range("I40").Value = 70052
$ws.Range("J40").Value = 70053
$ws.Range("K40").Value = 70052
$ws.Range("L40").Value = 70053
$ws.Range("M40").Value = -69892
$ws.Range("N40").Value = -70373
$ws.Range("H58").Value = 25003892
$ws.Range("I58").Value = 2904.4546
$ws.Range("J58").Value = 55560656
$ws.Range("K58").Value = 2904.4546
$ws.Range("L58").Value = 55560656
$ws.Range("M58").Value = -2701.4546
$ws.Range("N58").Value = -55561062
$ws.Range("H99").Value = 1987.5
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H126").Value = 1987.5
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H136").Value = 25003892
$ws.Range("I136").Value = 2904.4546
$ws.Range("J136").Value = 55560656
$ws.Range("K136").Value = 8713.363799999999
$ws.Range("L136").Value = 166681968
$ws.Range("M136").Value = -6163.363799999999
$ws.Range("N136").Value = -166687068

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 10763.7
$ws.Range("I6").Value = 107.4
$ws.Range("J6").Value = 21420
$ws.Range("K6").Value = 322.2
$ws.Range("L6").Value = 64260
$ws.Range("M6").Value = -209.2
$ws.Range("N6").Value = -64486
$ws.Range("H11").Value = 8185.077
$ws.Range("I11").Value = 281.2
$ws.Range("J11").Value = 13125
$ws.Range("K11").Value = 843.5999999999999
$ws.Range("L11").Value = 39375
$ws.Range("M11").Value = -703.5999999999999
$ws.Range("N11").Value = -39655
$ws.Range("H17").Value = 730
$ws.Range("J17").Value = 800
$ws.Range("L17").Value = 2400
$ws.Range("N17").Value = -2738
$ws.Range("H21").Value = 1399
$ws.Range("I21").Value = 920
$ws.Range("J21").Value = 1878
$ws.Range("K21").Value = 2760
$ws.Range("L21").Value = 5634
$ws.Range("M21").Value = -2587
$ws.Range("N21").Value = -5980
$ws.Range("H29").Value = 1048.3334
$ws.Range("I29").Value = 37.5
$ws.Range("J29").Value = 1337.1428
$ws.Range("K29").Value = 112.5
$ws.Range("L29").Value = 4011.4284
$ws.Range("M29").Value = 164.5
$ws.Range("N29").Value = -4565.428400000001
$ws.Range("H107").Value = 895.9
$ws.Range("I107").Value = 584.57574
$ws.Range("J107").Value = 1276.4073
$ws.Range("K107").Value = 1753.72722
$ws.Range("L107").Value = 3829.2219
$ws.Range("M107").Value = 166.27278
$ws.Range("N107").Value = -7669.2219
$ws.Range("H132").Value = 1510.2667
$ws.Range("I132").Value = 1252.1482
$ws.Range("J132").Value = 3833.3333
$ws.Range("K132").Value = 11269.3338
$ws.Range("L132").Value = 34499.9997
$ws.Range("M132").Value = -8739.3338
$ws.Range("N132").Value = -39559.9997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 29500
$ws.Range("J141").Value = 29500
$ws.Range("L141").Value = 29500
$ws.Range("N141").Value = -39860

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4024.5
$ws.Range("I132").Value = 2488.2222
$ws.Range("J132").Value = 5999.7144
$ws.Range("K132").Value = 7464.6666
$ws.Range("L132").Value = 17999.1432
$ws.Range("M132").Value = -4934.6666
$ws.Range("N132").Value = -23059.1432
$ws.Range("H136").Value = 3683.6667
$ws.Range("I136").Value = 3723.6924
$ws.Range("J136").Value = 3636.3635
$ws.Range("K136").Value = 11171.0772
$ws.Range("L136").Value = 10909.0905
$ws.Range("M136").Value = -8621.0772
$ws.Range("N136").Value = -16009.0905

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1607.2433
$ws.Range("I136").Value = 750.4828
$ws.Range("J136").Value = 4713
$ws.Range("K136").Value = 2251.4484
$ws.Range("L136").Value = 14139
$ws.Range("M136").Value = 298.5515999999998
$ws.Range("N136").Value = -19239

Write-Host "Applied all Atomos_Profits updates"